$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting the existing rows 5-7 down to 6-8
$ws.Rows("5:5").Insert()

# Copy the date style (s="2") from the row above so the new date cell keeps
# the same number format as the rest of column D
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the newly inserted row 5 with the new record's data
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"
$ws.Range("D5").Value = 44658
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 100112052
$ws.Range("G5").Value = "Albahaca"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 180
$ws.Range("K5").Value = 2500
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 2778
$ws.Range("N5").Value = "$/docena de matas"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 463
$ws.Range("Q5").Value = 6
$ws.Range("R5").Value = "Hortaliza"
